$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (row 2)
$ws.Range("G2").Value = 2

# Update cell values (row 3)
$ws.Range("D3").Value = 7
$ws.Range("G3").Value = 2

# Update cell values (row 4) - add new F4, G4 cells and update D4, E4, H4
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 4

# Update cell values (row 5)
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 1

# Update cell values (row 6)
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 3

# Update the selection/active cell
$ws.Range("E5").Select()
